$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-0.0413, -0.0413, -0.05119999870657921, 0.1909999996423721, -0.05119999870657921)
    3 = @(-0.0472, -0.0472, -0.0512, 0.1248999983072281, -0.05119999870657921)
    4 = @(-0.0508, -0.0508, -0.0512, 0.09380000084638596, -0.05119999870657921)
    5 = @(0.0522, 0.0522, 0.0644, 0.2610999941825867, -0.05119999870657921)
    6 = @(0.06850000000000001, 0.06850000000000001, 0.0441, 0.3282999992370605, -0.05119999870657921)
    7 = @(0.06610000000000001, 0.06610000000000001, 0.05990000069141388, 0.221000000834465, -0.05119999870657921)
    8 = @(0.08790000000000001, 0.08790000000000001, 0.08019999999999999, 0.221000000834465, -0.05119999870657921)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
